$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 162, shifting existing rows 162:224 down to 163:225
$ws.Rows.Item(162).EntireRow.Insert()

# Fill the new row 162 with its data
$ws.Cells.Item(162, 1).Value = 7
$ws.Cells.Item(162, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(162, 3).Value = "Ñuble"
$ws.Cells.Item(162, 4).Value = 44704
$ws.Cells.Item(162, 5).Value = 16
$ws.Cells.Item(162, 6).Value = 100112032
$ws.Cells.Item(162, 7).Value = "Zapallo italiano"
$ws.Cells.Item(162, 8).Value = "Sin especificar"
$ws.Cells.Item(162, 9).Value = "Primera"
$ws.Cells.Item(162, 10).Value = 60
$ws.Cells.Item(162, 11).Value = 17000
$ws.Cells.Item(162, 12).Value = 18000
$ws.Cells.Item(162, 13).Value = 17500
$ws.Cells.Item(162, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(162, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(162, 16).Value = 350
$ws.Cells.Item(162, 17).Value = 50
$ws.Cells.Item(162, 18).Value = "Hortaliza"
